# Auto-generated edit script: updates FFXIV Typhon market-profit values
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets per the scheduled-runner diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 163
$ws.Range("I4").Value = 163
$ws.Range("K4").Value = 163
$ws.Range("M4").Value = -49

$ws.Range("H5").Value = 242.78572
$ws.Range("I5").Value = 45.95
$ws.Range("J5").Value = 734.875
$ws.Range("K5").Value = 45.95
$ws.Range("L5").Value = 734.875
$ws.Range("M5").Value = 69.05
$ws.Range("N5").Value = -964.875

$ws.Range("H51").Value = 4810.4116
$ws.Range("J51").Value = 4187.7
$ws.Range("L51").Value = 4187.7
$ws.Range("N51").Value = -5155.7

$ws.Range("H62").Value = 2500
$ws.Range("I62").Value = 2500
$ws.Range("K62").Value = 2500
$ws.Range("M62").Value = -1876

$ws.Range("H65").Value = 2500
$ws.Range("I65").Value = 2500
$ws.Range("K65").Value = 12500
$ws.Range("M65").Value = -9380

$ws.Range("H86").Value = 16291.286
$ws.Range("I86").Value = 2584.3333
$ws.Range("J86").Value = 26571.5
$ws.Range("K86").Value = 2584.3333
$ws.Range("L86").Value = 26571.5
$ws.Range("M86").Value = -1461.3333
$ws.Range("N86").Value = -28817.5

$ws.Range("H89").Value = 16291.286
$ws.Range("I89").Value = 2584.3333
$ws.Range("J89").Value = 26571.5
$ws.Range("K89").Value = 12921.6665
$ws.Range("L89").Value = 132857.5
$ws.Range("M89").Value = -7305.666499999999
$ws.Range("N89").Value = -144089.5

$ws.Range("H92").Value = 504.73685
$ws.Range("J92").Value = 596.6667
$ws.Range("L92").Value = 596.6667
$ws.Range("N92").Value = -3092.6667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14713.338
$ws.Range("I32").Value = 15550.042
$ws.Range("K32").Value = 15550.042
$ws.Range("M32").Value = -15263.042

$ws.Range("H61").Value = 3960.111
$ws.Range("I61").Value = 3929.2856
$ws.Range("J61").Value = 3993.3076
$ws.Range("K61").Value = 3929.2856
$ws.Range("L61").Value = 3993.3076
$ws.Range("M61").Value = -3717.2856
$ws.Range("N61").Value = -4417.3076

$ws.Range("H74").Value = 1364.5
$ws.Range("I74").Value = 948.1
$ws.Range("J74").Value = 1885
$ws.Range("K74").Value = 948.1
$ws.Range("L74").Value = 1885
$ws.Range("M74").Value = -74.10000000000002
$ws.Range("N74").Value = -3633

$ws.Range("H77").Value = 1364.5
$ws.Range("I77").Value = 948.1
$ws.Range("J77").Value = 1885
$ws.Range("K77").Value = 4740.5
$ws.Range("L77").Value = 9425
$ws.Range("M77").Value = -372.5
$ws.Range("N77").Value = -18161

$ws.Range("H97").Value = 1895
$ws.Range("I97").Value = 1895
$ws.Range("K97").Value = 1895
$ws.Range("M97").Value = -1399

$ws.Range("H136").Value = 3960.111
$ws.Range("I136").Value = 3929.2856
$ws.Range("J136").Value = 3993.3076
$ws.Range("K136").Value = 11787.8568
$ws.Range("L136").Value = 11979.9228
$ws.Range("M136").Value = -9237.856800000001
$ws.Range("N136").Value = -17079.9228

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 352.8
$ws.Range("I22").Value = 352.8
$ws.Range("K22").Value = 352.8
$ws.Range("M22").Value = -179.8

$ws.Range("H99").Value = 2666.3333
$ws.Range("I99").Value = 2500
$ws.Range("J99").Value = 2999
$ws.Range("K99").Value = 2500
$ws.Range("L99").Value = 2999
$ws.Range("M99").Value = -1002
$ws.Range("N99").Value = -5995

$ws.Range("H105").Value = 2382531.5
$ws.Range("I105").Value = 1440
$ws.Range("K105").Value = 1440
$ws.Range("M105").Value = 307

$ws.Range("H134").Value = 35245.418
$ws.Range("I134").Value = 41641
$ws.Range("J134").Value = 1988.4
$ws.Range("K134").Value = 124923
$ws.Range("L134").Value = 5965.200000000001
$ws.Range("M134").Value = -122388
$ws.Range("N134").Value = -11035.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10650.082
$ws.Range("I31").Value = 24871.354
$ws.Range("J31").Value = 3095.0312
$ws.Range("K31").Value = 24871.354
$ws.Range("L31").Value = 3095.0312
$ws.Range("M31").Value = -24576.354
$ws.Range("N31").Value = -3685.0312

$ws.Range("H34").Value = 10650.082
$ws.Range("I34").Value = 24871.354
$ws.Range("J34").Value = 3095.0312
$ws.Range("K34").Value = 24871.354
$ws.Range("L34").Value = 3095.0312
$ws.Range("M34").Value = -24669.354
$ws.Range("N34").Value = -3499.0312

$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()  # cell removed entirely in target

$ws.Range("H58").Value = 11834.913
$ws.Range("I58").Value = 1013.3158
$ws.Range("J58").Value = 63237.5
$ws.Range("K58").Value = 1013.3158
$ws.Range("L58").Value = 63237.5
$ws.Range("M58").Value = -810.3158
$ws.Range("N58").Value = -63643.5

$ws.Range("H136").Value = 11834.913
$ws.Range("I136").Value = 1013.3158
$ws.Range("J136").Value = 63237.5
$ws.Range("K136").Value = 3039.9474
$ws.Range("L136").Value = 189712.5
$ws.Range("M136").Value = -489.9474
$ws.Range("N136").Value = -194812.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 207.5
$ws.Range("I14").Value = 207.5
$ws.Range("K14").Value = 622.5
$ws.Range("M14").Value = -449.5

$ws.Range("H23").Value = 755.58826
$ws.Range("I23").Value = 479.25
$ws.Range("K23").Value = 1437.75
$ws.Range("M23").Value = -1202.75

$ws.Range("H68").Value = 4788.893
$ws.Range("I68").Value = 645.0769
$ws.Range("J68").Value = 8380.200000000001
$ws.Range("K68").Value = 1935.2307
$ws.Range("L68").Value = 25140.6
$ws.Range("M68").Value = -1124.2307
$ws.Range("N68").Value = -26762.6

$ws.Range("H71").Value = 4788.893
$ws.Range("I71").Value = 645.0769
$ws.Range("J71").Value = 8380.200000000001
$ws.Range("K71").Value = 5805.6921
$ws.Range("L71").Value = 75421.8
$ws.Range("M71").Value = -1749.6921
$ws.Range("N71").Value = -83533.8

$ws.Range("H92").Value = 8675
$ws.Range("I92").Value = 700
$ws.Range("K92").Value = 2100
$ws.Range("M92").Value = -852

$ws.Range("H107").Value = 4630.885
$ws.Range("I107").Value = 25491.5
$ws.Range("K107").Value = 76474.5
$ws.Range("M107").Value = -74554.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12503500
$ws.Range("J70").Value = 20836500
$ws.Range("L70").Value = 20836500
$ws.Range("N70").Value = -20837040

$ws.Range("H73").Value = 12503500
$ws.Range("J73").Value = 20836500
$ws.Range("L73").Value = 20836500
$ws.Range("N73").Value = -20838372

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2740
$ws.Range("I22").Value = 3600.3333
$ws.Range("K22").Value = 3600.3333
$ws.Range("M22").Value = -3305.3333

$ws.Range("H27").Value = 2740
$ws.Range("I27").Value = 3600.3333
$ws.Range("K27").Value = 3600.3333
$ws.Range("M27").Value = -3493.3333

$ws.Range("H93").Value = 2548.6155
$ws.Range("J93").Value = 1966.6666
$ws.Range("L93").Value = 1966.6666
$ws.Range("N93").Value = -4462.6666

$ws.Range("H109").Value = 31992.5
$ws.Range("J109").Value = 31992.5
$ws.Range("L109").Value = 31992.5
$ws.Range("N109").Value = -34766.5

$ws.Range("H122").Value = 3556.2942
$ws.Range("I122").Value = 3250.25
$ws.Range("J122").Value = 3828.3333
$ws.Range("K122").Value = 9750.75
$ws.Range("L122").Value = 11484.9999
$ws.Range("M122").Value = -7300.75
$ws.Range("N122").Value = -16384.9999

$ws.Range("H136").Value = 17001.734
$ws.Range("I136").Value = 27608.21
$ws.Range("J136").Value = 3566.8667
$ws.Range("K136").Value = 82824.63
$ws.Range("L136").Value = 10700.6001
$ws.Range("M136").Value = -80274.63
$ws.Range("N136").Value = -15800.6001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 3862789.8
$ws.Range("I113").Value = 2400
$ws.Range("J113").Value = 13513764
$ws.Range("K113").Value = 7200
$ws.Range("L113").Value = 40541292
$ws.Range("M113").Value = -5030
$ws.Range("N113").Value = -40545632

$ws.Range("H136").Value = 1408.3704
$ws.Range("I136").Value = 944.4286
$ws.Range("J136").Value = 1908
$ws.Range("K136").Value = 2833.2858
$ws.Range("L136").Value = 5724
$ws.Range("M136").Value = -283.2857999999997
$ws.Range("N136").Value = -10824
